$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EPBDS-12787: the shared-string cell B8 held the literal (quote-prefixed)
# text "= addAll(null, null); "Hello";" -- rename the function reference to
# "addAll1" so the rule exercises the "no bruteforce implementation" path.
# A leading apostrophe keeps Excel from parsing the '=' as a formula,
# reproducing the quotePrefix cell style seen in the target file.
$ws.Range("B8").Value = "'= addAll1(null, null); ""Hello"";"

# The author's selection ended up on J8 when the workbook was last saved.
$ws.Range("J8").Select() | Out-Null
